$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target text "2014-04-21" looks like an ISO date, so assigning it
# directly to Range.Value/.Formula would make Excel auto-convert the cell
# to a date serial number (and bump its number format/style). To keep the
# cells as plain text (matching the original "Date" column layout), build
# the literal string in a scratch cell via a formula (whose cached result
# is a text value), then copy/paste-special *values only* into each target
# cell so no new style gets attached to the destination cells.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '="2014-04-21"'
$scratch.Copy()

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").PasteSpecial(-4163)  # xlPasteValues
}

$scratch.Clear()
$excel.CutCopyMode = $false
